$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.700.67"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "1.645.87"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'213.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'23.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.879.80"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.627.52"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D17").Value = "27.697.36"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "'231.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'10.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.38%  "
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").Value = "'150.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "1.439.34"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.570"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'0.886"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.42%  "
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'67.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").Value = "'2.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.788.64"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.05%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'85.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0988"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
